# Training Dashboard update — adding new progress as of date 04 Nov 2025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "STACKER CRANE AISLE MONTHLY PREVENTIVE MAINTENANCE (SOPs)" row (old row 3)
# is removed entirely; every row below it shifts up one.
$ws.Rows.Item(3).Delete()

# --- Row 3 (was row 4: "Endangered by Electricity A safety Training (SOPs)") ---
$ws.Range("A3").Value = 1
$ws.Range("H3").Value = -90

# Re-stamp the "LAST UPDATE" text cell with the new date without Excel's
# automatic text->date conversion: build it as a formula returning the
# literal string, then flatten the formula down to a plain value in place.
$ws.Range("I3").Formula = "=""04-Nov-2025"""
$ws.Range("I3").Copy()
$ws.Range("I3").PasteSpecial(-4163)

# --- Row 4 (was row 5: "Material request Procedure (Other Trainings)") ---
$ws.Range("A4").Value = 2
$ws.Range("H4").Value = 377

$ws.Range("I4").Formula = "=""04-Nov-2025"""
$ws.Range("I4").Copy()
$ws.Range("I4").PasteSpecial(-4163)

$excel.CutCopyMode = 0

# Column widths: B (TRAININGS) narrower, D (CODE) much narrower.
# ColumnWidth goes through Excel's MDW-based character<->internal-width
# conversion (adds ~5/6 of a character here), so back the request off by
# that fixed offset to land on the exact target width after round-trip.
$ws.Columns.Item(2).ColumnWidth = 52 - 5/6
$ws.Columns.Item(4).ColumnWidth = 6 - 5/6
